# NM2_2_KmCalc_5gDry — "Km calculation NM2-2 Day1"
#
# Adds the Day-1 (2017-08-22) measurement row to both the CH4 and CO2
# sheets, re-using row 4's number formatting, and leaves a review comment
# on each sheet's new date cell (C5).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "CH4" (first sheet)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("CH4")

$ws1.Range("C5").Value = 42969
$ws1.Range("D5").Formula = "=C5-C4"
$ws1.Range("E5").Value = 2734
$ws1.Range("F5").Value = 7077.67
$ws1.Range("G5").Value = 2633.64
$ws1.Range("H5").Value = 6816.12
$ws1.Range("I5").Value = 25.01
$ws1.Range("J5").Value = 64.49
$ws1.Range("K5").Value = 24.01
$ws1.Range("L5").Value = 61.9
$ws1.Range("M5").Value = 3
$ws1.Range("N5").Value = 7.6
$ws1.Range("O5").Value = 3.42
$ws1.Range("P5").Value = 8.6

# Match row 4's styling (date format on C, plain numeric elsewhere) by
# copying its formats down onto the new row instead of hand-picking
# style ids.
$ws1.Range("C4:P4").Copy()
$ws1.Range("C5:P5").PasteSpecial(-4122) | Out-Null

$ws1.Range("C5").AddComment("Author:`nethylene injections yesterday affected the methanizer, methane readings seem to be unaffected") | Out-Null

$ws1.Activate()
$ws1.Range("L10").Select()

# ---------------------------------------------------------------------
# Sheet "CO2" (second sheet)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("CO2")

$ws2.Range("C5").Value = 42969
$ws2.Range("D5").Formula = "=C5-C4"
$ws2.Range("E5").Value = 22.8
$ws2.Range("F5").Value = 0.079
$ws2.Range("G5").Value = 21.05
$ws2.Range("H5").Value = 0.0729
$ws2.Range("I5").Value = 21.47
$ws2.Range("J5").Value = 0.0743
$ws2.Range("K5").Value = 20.05
$ws2.Range("L5").Value = 0.0694
$ws2.Range("M5").Value = 19.86
$ws2.Range("N5").Value = 0.0688
$ws2.Range("O5").Value = 19.74
$ws2.Range("P5").Value = 0.0683

$ws2.Range("C4:P4").Copy()
$ws2.Range("C5:P5").PasteSpecial(-4122) | Out-Null

$ws2.Range("C5").AddComment("Author:`nDue to ethylene injections yesterday CO2 sensitivity has dropped at about 50%") | Out-Null

$ws2.Activate()
$ws2.Range("M10").Select()

Write-Host "NM2-2 Day1 rows + comments added"
